# Update "想去人数" (number of people interested) figures for the latest
# generated data pull, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1788
    $ws.Range("F4").Value = 549
    $ws.Range("F5").Value = 1144
    $ws.Range("F6").Value = 6029
}
